$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Requirements")

# Insert two new rows at the top of the Requirements sheet.
$ws1.Rows.Item(1).Insert()
$ws1.Rows.Item(1).Insert()

# Fill in the new header rows with group/author info.
$ws1.Range("A1").Value = "Group 22"
$ws1.Range("B1").Value = "עמאד טאהא emad.taha@e.braude.ac.il עיסא שבלי esa.shibli@e.braude.ac.il"
$ws1.Range("A2").Value = "Group 22"
$ws1.Range("B2").Value = "אחמד שחאדה ahmad.shhade@e.braude.ac.il יאמן אבו אחמד ואוי yamen.abu.ahmad.wawi@e.braude.ac.il"

# Switch the active sheet/selection to Requirements.
$ws1.Activate() | Out-Null
$ws1.Range("C14").Select() | Out-Null
